$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 1129.0834
$ws.Range("I18").Value = 1155.4
$ws.Range("J18").Value = 997.5
$ws.Range("K18").Value = 1155.4
$ws.Range("L18").Value = 997.5
$ws.Range("M18").Value = -871.4000000000001
$ws.Range("N18").Value = -1565.5

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 7995.222
$ws.Range("I32").Value = 11185.333
$ws.Range("J32").Value = 1615
$ws.Range("K32").Value = 11185.333
$ws.Range("L32").Value = 1615
$ws.Range("M32").Value = -10859.333
$ws.Range("N32").Value = -2267

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 275.3846
$ws.Range("I33").Value = 280.9091
$ws.Range("K33").Value = 280.9091
$ws.Range("M33").Value = -51.90910000000002

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 4001.5
$ws.Range("I40").Value = 3001
$ws.Range("K40").Value = 3001
$ws.Range("M40").Value = -2826

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H98").Value = 7302.8076
$ws.Range("I98").Value = 8694.380999999999
$ws.Range("J98").Value = 1458.2
$ws.Range("K98").Value = 8694.380999999999
$ws.Range("L98").Value = 1458.2
$ws.Range("M98").Value = -7196.380999999999
$ws.Range("N98").Value = -4454.2

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H100").Value = 8096.5
$ws.Range("J100").Value = 9195.549999999999
$ws.Range("L100").Value = 9195.549999999999
$ws.Range("N100").Value = -10277.55

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H106").Value = 3214.7
$ws.Range("I106").Value = 3041.1667
$ws.Range("K106").Value = 3041.1667
$ws.Range("M106").Value = -2410.1667

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H122").Value = 7302.8076
$ws.Range("I122").Value = 8694.380999999999
$ws.Range("J122").Value = 1458.2
$ws.Range("K122").Value = 26083.143
$ws.Range("L122").Value = 4374.6
$ws.Range("M122").Value = -23633.143
$ws.Range("N122").Value = -9274.6

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H131").Value = 2390.4243
$ws.Range("I131").Value = 1408.75
$ws.Range("J131").Value = 5008.222
$ws.Range("K131").Value = 4226.25
$ws.Range("L131").Value = 15024.666
$ws.Range("M131").Value = 813.75
$ws.Range("N131").Value = -25104.666

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 187.75
$ws.Range("I5").Value = 152.54546
$ws.Range("J5").Value = 575
$ws.Range("K5").Value = 152.54546
$ws.Range("L5").Value = 575
$ws.Range("M5").Value = -40.54545999999999
$ws.Range("N5").Value = -799

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 5595.5
$ws.Range("I74").Value = 3400
$ws.Range("K74").Value = 3400
$ws.Range("M74").Value = -2526

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 5595.5
$ws.Range("I77").Value = 3400
$ws.Range("K77").Value = 17000
$ws.Range("M77").Value = -12632

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H102").Value = 12504557
$ws.Range("I102").Value = 4507.6665
$ws.Range("K102").Value = 4507.6665
$ws.Range("M102").Value = -2885.6665

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 187.75
$ws.Range("I4").Value = 152.54546
$ws.Range("J4").Value = 575
$ws.Range("K4").Value = 152.54546
$ws.Range("L4").Value = 575
$ws.Range("M4").Value = -37.54545999999999
$ws.Range("N4").Value = -805

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 3998
$ws.Range("I99").Value = 3998
$ws.Range("K99").Value = 3998
$ws.Range("M99").Value = -2500

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 32274030
$ws.Range("I105").Value = 55576500
$ws.Range("K105").Value = 55576500
$ws.Range("M105").Value = -55574753

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 147.25
$ws.Range("I7").Value = 207.46153
$ws.Range("J7").Value = 35.42857
$ws.Range("K7").Value = 207.46153
$ws.Range("L7").Value = 35.42857
$ws.Range("M7").Value = -94.46153000000001
$ws.Range("N7").Value = -261.42857

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H7").Value = 834.6667
$ws.Range("J7").Value = 1002
$ws.Range("L7").Value = 3006
$ws.Range("N7").Value = -3230

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H107").Value = 1162.5758
$ws.Range("J107").Value = 1755
$ws.Range("L107").Value = 5265
$ws.Range("N107").Value = -9105

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H36").Value = 5374.75
$ws.Range("I36").Value = 749.5
$ws.Range("J36").Value = 10000
$ws.Range("K36").Value = 749.5
$ws.Range("L36").Value = 10000
$ws.Range("M36").Value = -264.5
$ws.Range("N36").Value = -10970

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H107").Value = 1455.909
$ws.Range("I107").Value = 708.6
$ws.Range("J107").Value = 2078.6667
$ws.Range("K107").Value = 708.6
$ws.Range("L107").Value = 2078.6667
$ws.Range("M107").Value = 1211.4
$ws.Range("N107").Value = -5918.6667

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 1237.5
$ws.Range("I122").Value = 1233.3334
$ws.Range("J122").Value = 1250
$ws.Range("K122").Value = 3700.0002
$ws.Range("L122").Value = 3750
$ws.Range("M122").Value = -1250.0002
$ws.Range("N122").Value = -8650

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 3695.2163
$ws.Range("I22").Value = 2493.625
$ws.Range("J22").Value = 4610.7144
$ws.Range("K22").Value = 2493.625
$ws.Range("L22").Value = 4610.7144
$ws.Range("M22").Value = -2198.625
$ws.Range("N22").Value = -5200.7144

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H27").Value = 3695.2163
$ws.Range("I27").Value = 2493.625
$ws.Range("J27").Value = 4610.7144
$ws.Range("K27").Value = 2493.625
$ws.Range("L27").Value = 4610.7144
$ws.Range("M27").Value = -2386.625
$ws.Range("N27").Value = -4824.7144

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 22738960
$ws.Range("I40").Value = 26325692
$ws.Range("J40").Value = 22998.666
$ws.Range("K40").Value = 26325692
$ws.Range("L40").Value = 22998.666
$ws.Range("M40").Value = -26325556
$ws.Range("N40").Value = -23270.666

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 4030.7856
$ws.Range("J61").Value = 4417.3335
$ws.Range("L61").Value = 4417.3335
$ws.Range("N61").Value = -4821.3335

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H100").Value = 11908017
$ws.Range("I100").Value = 41669000
$ws.Range("K100").Value = 41669000
$ws.Range("M100").Value = -41668459

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H113").Value = 4030.7856
$ws.Range("J113").Value = 4417.3335
$ws.Range("L113").Value = 4417.3335
$ws.Range("N113").Value = -8757.333500000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 1499
$ws.Range("J96").Value = 1499.3334
$ws.Range("L96").Value = 1499.3334
$ws.Range("N96").Value = -4245.3334

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 3466.5
$ws.Range("I107").Value = 1999.6666
$ws.Range("K107").Value = 5998.9998
$ws.Range("M107").Value = -4078.9998

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 2372.0334
$ws.Range("I122").Value = 2335.5925
$ws.Range("J122").Value = 2700
$ws.Range("K122").Value = 7006.7775
$ws.Range("L122").Value = 8100
$ws.Range("M122").Value = -4556.7775
$ws.Range("N122").Value = -13000

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 32125.97
$ws.Range("I126").Value = 35584.16
$ws.Range("J126").Value = 5325
$ws.Range("K126").Value = 106752.48
$ws.Range("L126").Value = 15975
$ws.Range("M126").Value = -104282.48
$ws.Range("N126").Value = -20915

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 6690.1304
$ws.Range("I132").Value = 5227.5713
$ws.Range("J132").Value = 8965.223
$ws.Range("K132").Value = 15682.7139
$ws.Range("L132").Value = 26895.669
$ws.Range("M132").Value = -13152.7139
$ws.Range("N132").Value = -31955.669
